$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename header column B from "file_name" to "code" to prepare context
# data for survey rendering.
$ws.Range("B1").Value = "code"

# Update the window/view state to match the saved selection.
$ws.Range("G16").Select()
